$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Conditional count" column (T). Excel shifts every
# column to its right one place to the left (U->T, V->U, W->V, X stays
# blank, Y->X, Z->Y, AA->Z) and re-points the formulas automatically.
$ws.Range("T:T").Delete()

# --- Column V (was "Fun Substitute" / SUBSTITUTE, now "Replace" / REPLACE) ---
$ws.Range("V1").Value = "Replace"
$ws.Range("V2").Formula = '=REPLACE(S2,4,1,"f")'
$ws.Range("V3").Formula = '=REPLACE(S3,4,1,"f")'

# --- Column Y (was "Count" referencing the old Y/Quality column, now
#     referencing the new X/Quality column after the shift) ---
$ws.Range("Y2").Formula = "=COUNTIF(I:I,X2)"
$ws.Range("Y3").Formula = "=COUNTIF(I:I,X3)"

# --- Column Z (was the array formula "Index value" =INDEX(Z:Z,2,0),
#     now a regular (non-array) formula =INDEX(Y2:Y100,1,1)) ---
$ws.Range("Z2").Formula = "=INDEX(Y2:Y100,1,1)"
$ws.Range("Z3").Formula = "=INDEX(Y2:Y100,1,1)"

# Selection moves to Z4 in the edited workbook.
$ws.Range("Z4").Select()
